# Updated ConDA result on full training set
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ConDA")

# ------------------------------------------------------------------
# Add 3 new rows (54:56) that replicate the "using z as the input to
# the classifier instead of h" / GaussianBlur augmentation setting
# (rows 42:44), but now reporting the result on the *entire* / full
# training set rather than the toy training set, with updated numbers.
# ------------------------------------------------------------------

# Copy the formatting (styles, borders, merged look, wrap text, etc.)
# of the existing block of rows (42:44) down onto the new block
# (54:56) without duplicating style definitions.
$src = $ws.Range("B42:I44")
$src.Copy()
$dst = $ws.Range("B54:I56")
$dst.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Match row heights of the source block.
$ws.Rows.Item(54).RowHeight = $ws.Rows.Item(42).RowHeight
$ws.Rows.Item(55).RowHeight = $ws.Rows.Item(43).RowHeight
$ws.Rows.Item(56).RowHeight = $ws.Rows.Item(44).RowHeight

# Re-create the merged cells for the new block (B, C, E, I columns).
$ws.Range("B54:B56").Merge() | Out-Null
$ws.Range("C54:C56").Merge() | Out-Null
$ws.Range("E54:E56").Merge() | Out-Null
$ws.Range("I54:I56").Merge() | Out-Null

# Row 54 - header row of the new block.
$ws.Cells.Item(54,2).Value2 = $ws.Cells.Item(42,2).Value2               # Model
$ws.Cells.Item(54,3).Value2 = "entire training set"                    # Training Set Used
$ws.Cells.Item(54,4).Value2 = $ws.Cells.Item(42,4).Value2               # Setting
$ws.Cells.Item(54,5).Value2 = $ws.Cells.Item(42,5).Value2               # Seed
$ws.Cells.Item(54,6).Value2 = $ws.Cells.Item(42,6).Value2               # Acc @ Covid
$ws.Cells.Item(54,7).Value2 = $ws.Cells.Item(42,7).Value2               # Acc @ Climate
$ws.Cells.Item(54,8).Value2 = 0.793                                     # Acc @ Military
$ws.Cells.Item(54,9).Value2 = $ws.Cells.Item(42,9).Value2               # Notes

# Row 55.
$ws.Cells.Item(55,4).Value2 = $ws.Cells.Item(43,4).Value2
$ws.Cells.Item(55,6).Value2 = $ws.Cells.Item(43,6).Value2
$ws.Cells.Item(55,7).Value2 = 0.808
$ws.Cells.Item(55,8).Value2 = $ws.Cells.Item(43,8).Value2

# Row 56.
$ws.Cells.Item(56,4).Value2 = $ws.Cells.Item(44,4).Value2
$ws.Cells.Item(56,6).Value2 = 0.801
$ws.Cells.Item(56,7).Value2 = $ws.Cells.Item(44,7).Value2
$ws.Cells.Item(56,8).Value2 = $ws.Cells.Item(44,8).Value2

# ------------------------------------------------------------------
# Update the view: scroll back to the top of the sheet and leave the
# selection on a single cell below the newly-added rows.
# ------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F62").Select() | Out-Null
